{"js": "// Update the worksheet date and each \"two-digit \u00f7 one-digit\" answer cell\n// to the next day's generated problem set (commit \"Update master to\n// output generated at 9a8706d\").\n//\n// Every \"from\" string below is unique within the document, so a\n// search-and-replace keyed on the old text is safe and unambiguous.\nconst replacements = [\n  [\"2024-02-12 Monday\", \"2024-02-13 Tuesday\"],\n  [\"93\u00f73=31, 0\", \"93\u00f74=23, 1\"],\n  [\"65\u00f73=21, 2\", \"86\u00f75=17, 1\"],\n  [\"81\u00f78=10, 1\", \"95\u00f73=31, 2\"],\n  [\"56\u00f74=14, 0\", \"78\u00f74=19, 2\"],\n  [\"12\u00f73=4, 0\", \"51\u00f79=5, 6\"],\n  [\"98\u00f79=10, 8\", \"97\u00f76=16, 1\"],\n  [\"10\u00f78=1, 2\", \"18\u00f78=2, 2\"],\n  [\"15\u00f73=5, 0\", \"85\u00f76=14, 1\"],\n  [\"16\u00f79=1, 7\", \"53\u00f78=6, 5\"],\n  [\"75\u00f74=18, 3\", \"64\u00f79=7, 1\"],\n  [\"14\u00f72=7, 0\", \"57\u00f76=9, 3\"],\n  [\"43\u00f73=14, 1\", \"80\u00f78=10, 0\"],\n  [\"69\u00f74=17, 1\", \"29\u00f72=14, 1\"],\n  [\"96\u00f75=19, 1\", \"97\u00f78=12, 1\"],\n  [\"28\u00f78=3, 4\", \"37\u00f78=4, 5\"],\n  [\"87\u00f72=43, 1\", \"90\u00f74=22, 2\"],\n  [\"89\u00f79=9, 8\", \"50\u00f76=8, 2\"],\n  [\"80\u00f75=16, 0\", \"33\u00f73=11, 0\"],\n  [\"42\u00f79=4, 6\", \"76\u00f74=19, 0\"],\n  [\"85\u00f79=9, 4\", \"53\u00f73=17, 2\"],\n  [\"78\u00f72=39, 0\", \"80\u00f72=40, 0\"],\n  [\"56\u00f78=7, 0\", \"55\u00f77=7, 6\"],\n  [\"59\u00f76=9, 5\", \"95\u00f75=19, 0\"],\n  [\"70\u00f73=23, 1\", \"33\u00f74=8, 1\"],\n  [\"87\u00f76=14, 3\", \"98\u00f75=19, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and each \"two-digit \u00f7 one-digit\" answer cell\n# to the next day's generated problem set (commit \"Update master to\n# output generated at 9a8706d\").\n#\n# Every \"from\" string below is unique within the document, so a\n# Find/Replace keyed on the old text is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-02-12 Monday\", \"2024-02-13 Tuesday\"),\n    @(\"93\u00f73=31, 0\", \"93\u00f74=23, 1\"),\n    @(\"65\u00f73=21, 2\", \"86\u00f75=17, 1\"),\n    @(\"81\u00f78=10, 1\", \"95\u00f73=31, 2\"),\n    @(\"56\u00f74=14, 0\", \"78\u00f74=19, 2\"),\n    @(\"12\u00f73=4, 0\", \"51\u00f79=5, 6\"),\n    @(\"98\u00f79=10, 8\", \"97\u00f76=16, 1\"),\n    @(\"10\u00f78=1, 2\", \"18\u00f78=2, 2\"),\n    @(\"15\u00f73=5, 0\", \"85\u00f76=14, 1\"),\n    @(\"16\u00f79=1, 7\", \"53\u00f78=6, 5\"),\n    @(\"75\u00f74=18, 3\", \"64\u00f79=7, 1\"),\n    @(\"14\u00f72=7, 0\", \"57\u00f76=9, 3\"),\n    @(\"43\u00f73=14, 1\", \"80\u00f78=10, 0\"),\n    @(\"69\u00f74=17, 1\", \"29\u00f72=14, 1\"),\n    @(\"96\u00f75=19, 1\", \"97\u00f78=12, 1\"),\n    @(\"28\u00f78=3, 4\", \"37\u00f78=4, 5\"),\n    @(\"87\u00f72=43, 1\", \"90\u00f74=22, 2\"),\n    @(\"89\u00f79=9, 8\", \"50\u00f76=8, 2\"),\n    @(\"80\u00f75=16, 0\", \"33\u00f73=11, 0\"),\n    @(\"42\u00f79=4, 6\", \"76\u00f74=19, 0\"),\n    @(\"85\u00f79=9, 4\", \"53\u00f73=17, 2\"),\n    @(\"78\u00f72=39, 0\", \"80\u00f72=40, 0\"),\n    @(\"56\u00f78=7, 0\", \"55\u00f77=7, 6\"),\n    @(\"59\u00f76=9, 5\", \"95\u00f75=19, 0\"),\n    @(\"70\u00f73=23, 1\", \"33\u00f74=8, 1\"),\n    @(\"87\u00f76=14, 3\", \"98\u00f75=19, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
